$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "Deposit" column (D) so it sits after "Founding Country" (was G),
# shifting Employee Size / Founding Year / Founding Country one column left.
$ws.Columns("D").Cut()
$ws.Columns("H").Insert()

# Deposit value correction for the NUS row (now in the relocated column G).
$ws.Range("G7").Value = 3000000000

# Net Income for the last row becomes the text "c" instead of 0.
$ws.Range("C8").Value = "c"
